$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows below have D (Price) and/or E (Volume(1h)) text values updated.
# Values are stored as text in the sheet (e.g. "297.58", "1.78%"), so each
# write forces a text number-format first and resets the style afterwards
# to avoid leaving a stray "Text" number format applied to the cell.
$updates = @(
    @{ Row = 2;  D = "297.58";         E = "1.78%" }
    @{ Row = 3;  D = "41.91";          E = "3.84%" }
    @{ Row = 4;  D = "5.000";          E = "-0.23%" }
    @{ Row = 5;  D = "0.07517";        E = "2.72%" }
    @{ Row = 6;  D = "1.585";          E = "4.14%" }
    @{ Row = 7;  D = "0.9268";         E = "0.14%" }
    @{ Row = 9;  D = "0.1199";         E = "-0.59%" }
    @{ Row = 10; D = "0.1832";         E = "5.78%" }
    @{ Row = 11; D = "0.08914";        E = "3.69%" }
    @{ Row = 12; D = "0.04084";        E = "-5.26%" }
    @{ Row = 13; E = "-0.72%" }
    @{ Row = 14; D = "0.001288";       E = "1.45%" }
    @{ Row = 15; D = "0.006006";       E = "0.82%" }
    @{ Row = 16; D = "3.357";          E = "0.57%" }
    @{ Row = 17; D = "4.385";          E = "2.08%" }
    @{ Row = 18; D = "0.3314";         E = "0.77%" }
    @{ Row = 19; D = "8.088";          E = "4.52%" }
    @{ Row = 20; E = "-0.01%" }
    @{ Row = 21; E = "18.21%" }
    @{ Row = 22; D = "0.04100";        E = "4.40%" }
    @{ Row = 23; E = "0.31%" }
    @{ Row = 24; D = "0.003888" }
    @{ Row = 25; E = "-3.91%" }
    @{ Row = 38; D = "0.02409";        E = "4.96%" }
    @{ Row = 39; D = "0.05233";        E = "5.14%" }
    @{ Row = 40; D = "0.006306";       E = "16.48%" }
    @{ Row = 41; D = "0.007824";       E = "1.71%" }
    @{ Row = 42; D = "0.1327";         E = "3.30%" }
    @{ Row = 43; D = "0.007391";       E = "0.88%" }
    @{ Row = 44; D = "0.007692";       E = "-2.86%" }
    @{ Row = 45; D = "0.2962";         E = "-6.89%" }
    @{ Row = 46; D = "0.00006467";     E = "2.27%" }
    @{ Row = 47; D = "0.00000000751";  E = "0.02%" }
    @{ Row = 48; D = "0.03133";        E = "53.44%" }
    @{ Row = 49; D = "0.004203";       E = "0.04%" }
    @{ Row = 50; D = "0.00002102";     E = "0.02%" }
    @{ Row = 51; D = "0.0002002";      E = "0.02%" }
)

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey("D")) {
        Set-TextCell $ws.Cells.Item($r, 4) $u.D
    }
    if ($u.ContainsKey("E")) {
        Set-TextCell $ws.Cells.Item($r, 5) $u.E
    }
}
